# This script swaps the data of three pairs of rows in the worksheet
# (13<->15, 17<->18, 22<->24), as described by the source diff.
# Row numbers / record positions stay where they are, but the observation
# data that lives in each row moves to the paired row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($row1, $row2) {
    # Only the columns whose content actually differs between the two rows
    # of a swapped pair (per the source diff):
    # A=1 Id, B=2 Taxonsorteringsordning, D=4 Rodlistade, E=5 TaxonId,
    # F=6 Artnamn, G=7 Vetenskapligt namn, H=8 Auktor, Q=17 Ost, R=18 Nord,
    # S=19 Noggrannhet, Z=26 Starttid, AB=28 Sluttid, AC=29 Publik kommentar
    $cols = @(1,2,4,5,6,7,8,17,18,19,26,28,29)

    foreach ($col in $cols) {
        $cell1 = $ws.Cells.Item($row1, $col)
        $cell2 = $ws.Cells.Item($row2, $col)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}

Swap-RowData 13 15
Swap-RowData 17 18
Swap-RowData 22 24
